$d = $word.ActiveDocument

# Unicode helpers (right single quote, en dash) so the literal text
# matches the source exactly regardless of script-file encoding.
$rsq = [char]0x2019
$dash = [char]0x2013

# -----------------------------------------------------------------
# 1) Insert a new "Start of week:" Heading1 paragraph (double
#    underline) right before the existing "Evaluation" Heading1
#    paragraph.
# -----------------------------------------------------------------
$evalIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd("`r`a") -eq "Evaluation") {
        $evalIdx = $i
        break
    }
}

$evalPara = $d.Paragraphs($evalIdx)
$evalPara.Range.InsertParagraphBefore()
$startPara = $d.Paragraphs($evalIdx)
$startPara.Style = "Heading1"
$startPara.Range.Text = "Start of week:"
$startPara.Range.Font.Underline = 3

# -----------------------------------------------------------------
# 2) Append the "End of week" section + per-person updates at the
#    very end of the document (after the last paragraph).
#    Build everything as plain/"Normal" + justified first (so the
#    double-underline heading formatting is not copied forward onto
#    later paragraphs), then retro-fit the heading paragraph's style
#    as the final step.
# -----------------------------------------------------------------
$idx = $d.Paragraphs.Count

# -- "End of week - Who did what and Progression:" heading ---------
$cur = $d.Paragraphs($idx)
$cur.Range.InsertParagraphAfter()
$idx = $idx + 1
$cur = $d.Paragraphs($idx)
$cur.Alignment = 3
$cur.Range.Text = "End of week " + $dash + " Who did what and Progression:"
$headingIdx = $idx

# -- blank line -------------------------------------------------------
$cur.Range.InsertParagraphAfter()
$idx = $idx + 1
$cur = $d.Paragraphs($idx)
$cur.Alignment = 3

# -- Chenlei paragraph --------------------------------------------------
$cur.Range.InsertParagraphAfter()
$idx = $idx + 1
$cur = $d.Paragraphs($idx)
$cur.Alignment = 3
$chenleiText = "Chenlei: Wasn" + $rsq + "t familiar with Cloud 9. Created account and tested small bit of code to see how it all worked. Tried using phpMyAdmin but having issues. Created form and PHP coding. Tried to connect them together, currently not working. Will consult with team on Monday to discuss resolution. "
$cur.Range.Text = $chenleiText
$boldEnd = $cur.Range.Start + ("Chenlei:").Length
$boldRng = $d.Range($cur.Range.Start, $boldEnd)
$boldRng.Font.Bold = 1

# -- blank line -------------------------------------------------------
$cur.Range.InsertParagraphAfter()
$idx = $idx + 1
$cur = $d.Paragraphs($idx)
$cur.Alignment = 3

# -- Keith paragraph ------------------------------------------------------
$cur.Range.InsertParagraphAfter()
$idx = $idx + 1
$cur = $d.Paragraphs($idx)
$cur.Alignment = 3
$keithBody = "Used Cloud 9 before. Jumped straight in. Created Customer registration in User-Reg-Login. Created MySQL database, a form and Form validation to avoid SQL injections. Created error page. Linked PHP to form. Thought it would be a good idea to have email verification. Created a create_password.php page where the link in the email would go"
$keithText = "Keith: " + $keithBody + ". Have started on it. "
$cur.Range.Text = $keithText
$boldEnd = $cur.Range.Start + ("Keith:").Length
$boldRng = $d.Range($cur.Range.Start, $boldEnd)
$boldRng.Font.Bold = 1
# re-home the (single, Word-managed) "_GoBack" bookmark at the point
# of this latest edit, matching Word's own behaviour.
$goPos = $cur.Range.Start + ("Keith: " + $keithBody).Length
$goRng = $d.Range($goPos, $goPos)
$d.Bookmarks.Add("_GoBack", $goRng)

# -- blank line -------------------------------------------------------
$cur.Range.InsertParagraphAfter()
$idx = $idx + 1
$cur = $d.Paragraphs($idx)
$cur.Alignment = 3

# -- Ali paragraph -------------------------------------------------------
$cur.Range.InsertParagraphAfter()
$idx = $idx + 1
$cur = $d.Paragraphs($idx)
$cur.Alignment = 3
$cur.Range.Text = "Ali: *Nothing provided. Not in class Thursday and Friday*"
$boldEnd = $cur.Range.Start + ("Ali:").Length
$boldRng = $d.Range($cur.Range.Start, $boldEnd)
$boldRng.Font.Bold = 1

# -- blank line -------------------------------------------------------
$cur.Range.InsertParagraphAfter()
$idx = $idx + 1
$cur = $d.Paragraphs($idx)
$cur.Alignment = 3

# -- Surendra paragraph -------------------------------------------------
$cur.Range.InsertParagraphAfter()
$idx = $idx + 1
$cur = $d.Paragraphs($idx)
$cur.Alignment = 3
$cur.Range.Text = "Surendra: [input your section here]"
$boldEnd = $cur.Range.Start + ("Surendra:").Length
$boldRng = $d.Range($cur.Range.Start, $boldEnd)
$boldRng.Font.Bold = 1

# -- five trailing blank lines -------------------------------------
for ($i = 0; $i -lt 5; $i++) {
    $cur.Range.InsertParagraphAfter()
    $idx = $idx + 1
    $cur = $d.Paragraphs($idx)
    $cur.Alignment = 3
}

# -- retro-fit the "End of week..." heading's style/underline ------
$headingPara = $d.Paragraphs($headingIdx)
$headingPara.Style = "Heading1"
$headingPara.Range.Font.Underline = 3

Write-Output "Done"
